# Apply the "Add files via upload" revision to the stickers workbook:
#  - drop the now-unused second sheet ("Лист2")
#  - add a "file_id" column header in B1 and a new "Ответ" column in C
#  - fill in bot reply text for the first couple of rows
#  - tidy the question phrases ("как дела?" -> "как дела", "знаешь?" -> "знаешь")
#  - restore the active-cell selection and widen column A a bit

$wb = $excel.ActiveWorkbook

# --- remove the empty second worksheet ("Лист2") ---------------------------
foreach ($sheet in @($wb.Worksheets)) {
    if ($sheet.Name -ne "stickers") {
        [void]$sheet.Delete()
    }
}

$ws = $wb.Worksheets.Item("stickers")

# --- header row: add file_id label + new "Ответ" column --------------------
$ws.Range("B1").Value = "file_id"
$ws.Range("C1").Value = "Ответ"

# --- new bot-reply text in column C for the first two sticker rows ---------
$ws.Range("C2").Value = "Приятно познакомиться с живым человеком, я - бот!"
$ws.Range("C3").Value = "жаль расстоваться, человек)"

# --- tidy up the trigger phrases --------------------------------------------
$ws.Range("A4").Value = "как дела"
$ws.Range("A6").Value = "знаешь"

# --- cosmetic touches: column width + selected cell -------------------------
$ws.Range("A1").ColumnWidth = 15.25
[void]$ws.Range("A4").Select()
